$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.186996666666667
$ws.Range("H2").Value = 6.56099
$ws.Range("I2").Value = 0.9253339937566305
$ws.Range("J2").Value = 0.9253339937566305
$ws.Range("M2").Value = 0.06166766666666667
$ws.Range("N2").Value = 0.185003
$ws.Range("O2").Value = 0.3189772891852935
$ws.Range("P2").Value = 0.3189772891852935
$ws.Range("Q2").Value = 0.1348669814411111
$ws.Range("R2").Value = 1.21380283297
$ws.Range("S2").Value = 0.2951605289194913
$ws.Range("T2").Value = 0.2951605289194913
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.186996666666667
$ws.Range("H3").Value = 6.56099
$ws.Range("I3").Value = 0.9253339937566305
$ws.Range("J3").Value = 0.9253339937566305
$ws.Range("O3").Value = 0.4045463009579509
$ws.Range("P3").Value = 0.4045463009579509
$ws.Range("Q3").Value = 0.1710464672977778
$ws.Range("R3").Value = 1.53941820568
$ws.Range("S3").Value = 0.3743404443248924
$ws.Range("T3").Value = 0.3743404443248925
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.186996666666667
$ws.Range("H4").Value = 6.56099
$ws.Range("I4").Value = 0.9253339937566305
$ws.Range("J4").Value = 0.9253339937566305
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05345100000000001
$ws.Range("N4").Value = 0.160353
$ws.Range("O4").Value = 0.2764764098567557
$ws.Range("P4").Value = 0.2764764098567557
$ws.Range("Q4").Value = 0.11689715883
$ws.Range("R4").Value = 1.05207442947
$ws.Range("S4").Value = 0.2558330205122468
$ws.Range("T4").Value = 0.2558330205122468
$ws.Range("I5").Value = 0.07466600624336955
$ws.Range("J5").Value = 0.07466600624336955
$ws.Range("M5").Value = 0.06166766666666667
$ws.Range("N5").Value = 0.185003
$ws.Range("O5").Value = 0.3189772891852935
$ws.Range("P5").Value = 0.3189772891852935
$ws.Range("Q5").Value = 0.01088253424844444
$ws.Range("R5").Value = 0.097942808236
$ws.Range("S5").Value = 0.02381676026580222
$ws.Range("T5").Value = 0.02381676026580222
$ws.Range("I6").Value = 0.07466600624336955
$ws.Range("J6").Value = 0.07466600624336955
$ws.Range("O6").Value = 0.4045463009579509
$ws.Range("P6").Value = 0.4045463009579509
$ws.Range("S6").Value = 0.03020585663305842
$ws.Range("T6").Value = 0.03020585663305842
$ws.Range("I7").Value = 0.07466600624336955
$ws.Range("J7").Value = 0.07466600624336955
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.05345100000000001
$ws.Range("N7").Value = 0.160353
$ws.Range("O7").Value = 0.2764764098567557
$ws.Range("P7").Value = 0.2764764098567557
$ws.Range("Q7").Value = 0.009432533604000001
$ws.Range("R7").Value = 0.08489280243600002
$ws.Range("S7").Value = 0.02064338934450892
$ws.Range("T7").Value = 0.02064338934450892
